# "Se avanza con Productos"
# - Clears the AutoFilter criteria on the Status/Tarea sheet (was filtering
#   column G ("Persona") to only show "Sebas"), which also unhides every
#   previously filtered-out row.
# - Marks the "Products" tasks (rows 33/34 - F33/F34) as finished instead of
#   pending, matching the status color used by the other "finished" cells.
# - Moves the sheet's active selection down to F36 (where the user was last
#   working), mirroring the scrolled view position.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Clear all filtering (un-hides the rows Excel had hidden for the "Sebas"
# filter on the Persona column, and drops the stored filter criteria).
$ws.ShowAllData()

# "Crear Pantalla Products" / "Crear Metodos de ABM para manejar productos"
# status: pending -> finished. Copy the formatting from an existing
# "finished" cell (F6) so the fill/alignment match, then set the text.
$ws.Range("F6").Copy()
$ws.Range("F33").PasteSpecial(-4122)
$ws.Range("F34").PasteSpecial(-4122)
$ws.Range("F33").Value = "finished"
$ws.Range("F34").Value = "finished"

# Move the selection to where work continued.
$ws.Activate()
$ws.Range("F36").Select()
